$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.167224666666667
$ws.Range("H2").Value = 9.501674000000001
$ws.Range("I2").Value = 0.007330706330134895
$ws.Range("J2").Value = 0.007330706330134894
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 244.1598398429089
$ws.Range("R2").Value = 2197.438558586181
$ws.Range("S2").Value = 0.001762167410657304
$ws.Range("T2").Value = 0.001762167410657304
$ws.Range("G3").Value = 3.167224666666667
$ws.Range("H3").Value = 9.501674000000001
$ws.Range("I3").Value = 0.007330706330134895
$ws.Range("J3").Value = 0.007330706330134894
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 321.7267998830543
$ws.Range("R3").Value = 2895.541198947488
$ws.Range("S3").Value = 0.002321989079996719
$ws.Range("T3").Value = 0.002321989079996718
$ws.Range("G4").Value = 3.167224666666667
$ws.Range("H4").Value = 9.501674000000001
$ws.Range("I4").Value = 0.007330706330134895
$ws.Range("J4").Value = 0.007330706330134894
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 449.8307505040037
$ws.Range("R4").Value = 4048.476754536033
$ws.Range("S4").Value = 0.003246549839480873
$ws.Range("T4").Value = 0.003246549839480873
$ws.Range("I5").Value = 0.9163769820903682
$ws.Range("J5").Value = 0.9163769820903681
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 30521.26863453216
$ws.Range("R5").Value = 274691.4177107894
$ws.Range("S5").Value = 0.2202802268968294
$ws.Range("T5").Value = 0.2202802268968294
$ws.Range("I6").Value = 0.9163769820903682
$ws.Range("J6").Value = 0.9163769820903681
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2902608902538085
$ws.Range("T6").Value = 0.2902608902538084
$ws.Range("I7").Value = 0.9163769820903682
$ws.Range("J7").Value = 0.9163769820903681
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.4058358649397304
$ws.Range("T7").Value = 0.4058358649397303
$ws.Range("I8").Value = 0.07629231157949697
$ws.Range("J8").Value = 0.07629231157949695
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 2541.026435600308
$ws.Range("R8").Value = 22869.23792040277
$ws.Range("S8").Value = 0.01833927306519569
$ws.Range("T8").Value = 0.01833927306519568
$ws.Range("I9").Value = 0.07629231157949697
$ws.Range("J9").Value = 0.07629231157949695
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.02416546324425458
$ws.Range("T9").Value = 0.02416546324425457
$ws.Range("I10").Value = 0.07629231157949697
$ws.Range("J10").Value = 0.07629231157949695
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.03378757527004671
$ws.Range("T10").Value = 0.0337875752700467